# Auto-generated COM script to apply the Inflation_Targets sheet population
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("HL_Inflation")
$ws2 = $wb.Worksheets.Item("Inflation_Targets")

# Row 2: country header row (bold)
$ws2.Range("B2").Value = 'India'
$ws2.Range("B2").Font.Bold = $true
$ws2.Range("C2").Value = 'Singapore'
$ws2.Range("C2").Font.Bold = $true
$ws2.Range("D2").Value = 'Canada'
$ws2.Range("D2").Font.Bold = $true
$ws2.Range("E2").Value = 'USA'
$ws2.Range("E2").Font.Bold = $true
$ws2.Range("F2").Value = 'Japan'
$ws2.Range("F2").Font.Bold = $true
$ws2.Range("G2").Value = 'Australia'
$ws2.Range("G2").Font.Bold = $true
$ws2.Range("H2").Value = 'China'
$ws2.Range("H2").Font.Bold = $true
$ws2.Range("I2").Value = 'Switzerland'
$ws2.Range("I2").Font.Bold = $true
$ws2.Range("J2").Value = 'UK'
$ws2.Range("J2").Font.Bold = $true
$ws2.Range("K2").Value = 'Spain'
$ws2.Range("K2").Font.Bold = $true
$ws2.Range("L2").Value = 'France'
$ws2.Range("L2").Font.Bold = $true
$ws2.Range("M2").Value = 'Germany'
$ws2.Range("M2").Font.Bold = $true
$ws2.Range("N2").Value = 'Italy'
$ws2.Range("N2").Font.Bold = $true

# Row 3: target rate values
$ws2.Range("A3").Value = 'Target Rate'
$ws2.Range("B3").Value = '0.02 - 0.06'
$ws2.Range("C3").Value = 0.02
$ws2.Range("D3").Value = 0.02
$ws2.Range("E3").Value = 0.02
$ws2.Range("F3").Value = 0.02
$ws2.Range("F3").Font.Color = 0
$ws2.Range("G3").Value = '0.02 - 0.03'
$ws2.Range("H3").Value = 0.02
$ws2.Range("I3").Value = 0.02
$ws2.Range("J3").Value = 0.02
$ws2.Range("K3").Value = 0.02
$ws2.Range("L3").Value = 0.02
$ws2.Range("M3").Value = 0.02
$ws2.Range("N3").Value = 0.02

# Row 4: since-dates / notes
$ws2.Range("B4").Value = 'since 31.03.2016'
$ws2.Range("E4").Value = 'since 31.12.2011'
$ws2.Range("F4").Value = 'since 22.01.2013'
$ws2.Range("H4").Value = '"reflects our interpretation of reasonable and stable" inflation'
$ws2.Range("J4").Value = 'since 31.12.2003'

# Row 5
$ws2.Range("B5").Value = 0.04
$ws2.Range("E5").Value = '0.017 - 0.02'

# Row 6
$ws2.Range("B6").Value = 'before 31.03.2016'
$ws2.Range("E6").Value = 'before 2012'

# Rows 22-34: country + description table
$ws2.Range("A22").Value = 'India'
$ws2.Range("A22").Font.Bold = $true
$ws2.Range("B22").Value = 'India: The Reserve Bank of India has changed its target inflation rate over time. From 1998-2015, the target was an inflation rate of 4% with a +/- 2% tolerance range. In 2016, the RBI adopted a new monetary policy framework, and the target was changed to a headline inflation rate of 4% (+/- 2%) until March 2021. In April 2021, the RBI announced a new target range of 2-6% for the next five years.'
$ws2.Range("A23").Value = 'Singapore'
$ws2.Range("A23").Font.Bold = $true
$ws2.Range("B23").Value = 'Singapore: The Monetary Authority of Singapore has adopted a policy of targeting medium-term price stability, and does not have a specific numerical target for inflation.'
$ws2.Range("A24").Value = 'Canada'
$ws2.Range("A24").Font.Bold = $true
$ws2.Range("B24").Value = 'Canada: The Bank of Canada has changed its target inflation rate over time. From 1991 to 2001, the target was a range of 1% to 3%. From 2002 to 2011, the target was a 2% inflation rate. In 2012, the Bank of Canada adopted a new policy of targeting a 2% inflation rate that is "flexible" in the short term.'
$ws2.Range("A25").Value = 'USA'
$ws2.Range("A25").Font.Bold = $true
$ws2.Range("B25").Value = 'USA: The Federal Reserve has changed its target inflation rate over time. From 1996 to 2011, the target was a 2% inflation rate. In August 2020, the Federal Reserve announced a new policy of "average inflation targeting," where it would aim to achieve inflation that averages 2% over time, rather than targeting a specific rate.'
$ws2.Range("A26").Value = 'Japan'
$ws2.Range("A26").Font.Bold = $true
$ws2.Range("B26").Value = 'Japan: The Bank of Japan has had a target inflation rate of 2% since 2013, as part of its "Quantitative and Qualitative Monetary Easing" policy. Prior to that, the Bank of Japan had various inflation targets, including a 1% target from 2010 to 2013.'
$ws2.Range("A27").Value = 'Australia'
$ws2.Range("A27").Font.Bold = $true
$ws2.Range("B27").Value = 'Australia: The Reserve Bank of Australia has changed its target inflation rate over time. From 1993 to 2016, the target was a 2% to 3% inflation rate. In 2016, the target was changed to a 2% inflation rate, with a "flexible medium-term inflation target."'
$ws2.Range("A28").Value = 'China'
$ws2.Range("A28").Font.Bold = $true
$ws2.Range("B28").Value = 'China: The People''s Bank of China does not have a specific inflation target, but aims for "reasonable and stable" inflation.'
$ws2.Range("A29").Value = 'Switzerland'
$ws2.Range("A29").Font.Bold = $true
$ws2.Range("B29").Value = 'Switzerland: The Swiss National Bank has had a target inflation rate of below 2% since 2000.'
$ws2.Range("A30").Value = 'UK'
$ws2.Range("A30").Font.Bold = $true
$ws2.Range("B30").Value = 'UK: The Bank of England has changed its target inflation rate over time. From 1992 to 2003, the target was a 2.5% inflation rate. From 2004 to 2013, the target was a 2% inflation rate. In 2013, the Bank of England adopted a "forward guidance" policy, where it would not consider raising interest rates until the unemployment rate fell below 7%. In 2015, the Bank of England returned to targeting a 2% inflation rate.'
$ws2.Range("A31").Value = 'Spain'
$ws2.Range("A31").Font.Bold = $true
$ws2.Range("B31").Value = 'Spain: The Bank of Spain does not have a specific inflation target, but aims for price stability in the eurozone.'
$ws2.Range("A32").Value = 'France'
$ws2.Range("A32").Font.Bold = $true
$ws2.Range("B32").Value = 'France: The Bank of France does not have a specific inflation target, but aims for price stability in the eurozone.'
$ws2.Range("A33").Value = 'Germany'
$ws2.Range("A33").Font.Bold = $true
$ws2.Range("B33").Value = 'Germany: The Bundesbank does not have a specific inflation target, but aims for price stability in the eurozone.'
$ws2.Range("A34").Value = 'Italy'
$ws2.Range("A34").Font.Bold = $true
$ws2.Range("B34").Value = 'Italy: The Bank of Italy does not have a specific inflation target, but aims for price stability in the eurozone.'

# Selections
$ws2.Activate()
$ws2.Range("H5").Select()
$ws1.Activate()
$ws1.Range("F15").Select()

